$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin market data (price, volume, name/link shifts) per commit

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '278.45'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '6.52%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '27.28'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '1.33%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '4.792'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '1.61%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.06293'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '1.24%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.924'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '2.91%'
$ws.Range('B7').NumberFormat = '@'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').NumberFormat = '@'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.270'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '2.85%'
$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').NumberFormat = '@'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8780'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '3.31%'
$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'FTXToken'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9429'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '3.25%'
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1462'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '4.08%'
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.05166'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '4.65%'
$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07277'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '2.63%'
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03124'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '0.29%'
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.09069'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '0.16%'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001555'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '1.73%'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0006287'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '1.76%'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.005879'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-1.42%'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.450'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '-0.01%'
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'BTSEToken'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.284'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '5.39%'
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'BitpandaEcosystemToken'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.3147'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '1.61%'
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'ProBitToken'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1312'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '0.15%'
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'MCDex'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.849'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-5.99%'
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'CoinExToken'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04335'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '2.21%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001181'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-0.01%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004284'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '5.18%'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-0.04%'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0001690'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '3.06%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04067'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '3.02%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006588'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '58.97%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1156'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '3.85%'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '1.53%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005113'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-0.95%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000750'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-0.05%'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '856.84%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.02251'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '-33.85%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002100'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.05%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0002000'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.05%'
